$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '42.920.05'
$ws.Range("E2").Value = "'" + '  -0.86%  '
$ws.Range("D3").Value = "'" + '2.237.63'
$ws.Range("E3").Value = "'" + '  -2.13%  '
$ws.Range("E4").Value = "'" + '  +0.32%  '
$ws.Range("D5").Value = "'" + '112.41'
$ws.Range("E5").Value = "'" + '  -1.29%  '
$ws.Range("D6").Value = "'" + '275.71'
$ws.Range("E6").Value = "'" + '  +3.57%  '
$ws.Range("D7").Value = "'" + '0.626'
$ws.Range("E7").Value = "'" + '  +0.06%  '
$ws.Range("D8").Value = "'" + '1.00'
$ws.Range("E8").Value = "'" + '  +0.28%  '
$ws.Range("D9").Value = "'" + '0.606'
$ws.Range("E9").Value = "'" + '  -0.68%  '
$ws.Range("D10").Value = "'" + '46.02'
$ws.Range("E10").Value = "'" + '  -2.82%  '
$ws.Range("D11").Value = "'" + '0.0927'
$ws.Range("E11").Value = "'" + '  -1.05%  '
$ws.Range("D12").Value = "'" + '9.02'
$ws.Range("E12").Value = "'" + '  -3.55%  '
$ws.Range("E13").Value = "'" + '  -2.97%  '
$ws.Range("D14").Value = "'" + '15.27'
$ws.Range("E14").Value = "'" + '  -1.78%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = "'" + '2.575.02'
$ws.Range("E15").Value = "'" + '  -2.01%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = "'" + '0.867'
$ws.Range("E16").Value = "'" + '  -0.21%  '
$ws.Range("D17").Value = "'" + '2.239.70'
$ws.Range("E17").Value = "'" + '  -1.58%  '
$ws.Range("D18").Value = "'" + '42.827.86'
$ws.Range("E18").Value = "'" + '  -1.12%  '
$ws.Range("E19").Value = "'" + '  -1.68%  '
$ws.Range("D20").Value = "'" + '6.74'
$ws.Range("E20").Value = "'" + '  -1.33%  '
$ws.Range("D21").Value = "'" + '71.91'
$ws.Range("E21").Value = "'" + '  -0.25%  '
$ws.Range("D22").Value = "'" + '2.33'
$ws.Range("E22").Value = "'" + '  -5.60%  '
$ws.Range("D23").Value = "'" + '2.99'
$ws.Range("E23").Value = "'" + '  +3.61%  '
$ws.Range("D24").Value = "'" + '230.97'
$ws.Range("E24").Value = "'" + '  -1.68%  '
$ws.Range("D25").Value = "'" + '9.24'
$ws.Range("E25").Value = "'" + '  -3.70%  '
$ws.Range("D26").Value = "'" + '12.14'
$ws.Range("E26").Value = "'" + '  +5.94%  '
$ws.Range("D28").Value = "'" + '40.29'
$ws.Range("E28").Value = "'" + '  -2.97%  '
$ws.Range("D29").Value = "'" + '2.24'
$ws.Range("E29").Value = "'" + '  -0.57%  '
$ws.Range("D30").Value = "'" + '3.26'
$ws.Range("E30").Value = "'" + '  -2.58%  '
$ws.Range("D31").Value = "'" + '173.94'
$ws.Range("E31").Value = "'" + '  +0.02%  '
$ws.Range("D32").Value = "'" + '21.08'
$ws.Range("E32").Value = "'" + '  -2.67%  '
$ws.Range("D33").Value = "'" + '0.0901'
$ws.Range("E33").Value = "'" + '  -0.83%  '
$ws.Range("D34").Value = "'" + '5.55'
$ws.Range("E34").Value = "'" + '  -2.92%  '
$ws.Range("D35").Value = "'" + '4.31'
$ws.Range("E35").Value = "'" + '  +7.71%  '
$ws.Range("D36").Value = "'" + '0.127'
$ws.Range("E36").Value = "'" + '  -0.51%  '
$ws.Range("D37").Value = "'" + '4.64'
$ws.Range("E37").Value = "'" + '  +0.23%  '
$ws.Range("D38").Value = "'" + '0.0371'
$ws.Range("E38").Value = "'" + '  +0.63%  '
$ws.Range("D39").Value = "'" + '0.106'
$ws.Range("E39").Value = "'" + '  +1.68%  '
$ws.Range("D40").Value = "'" + '2.56'
$ws.Range("E40").Value = "'" + '  -3.03%  '
$ws.Range("D41").Value = "'" + '70.83'
$ws.Range("E41").Value = "'" + '  -7.34%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").Value = "'" + '13.11'
$ws.Range("E42").Value = "'" + '  -8.56%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = "'" + '0.231'
$ws.Range("E43").Value = "'" + '  -3.83%  '
$ws.Range("E44").Value = "'" + '  +0.01%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = "'" + '1.33'
$ws.Range("E45").Value = "'" + '  -4.32%  '
$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D46").Value = "'" + '5.63'
$ws.Range("E46").Value = "'" + '  -8.70%  '
$ws.Range("D47").Value = "'" + '1.27'
$ws.Range("E47").Value = "'" + '  +0.76%  '
$ws.Range("D48").Value = "'" + '8.39'
$ws.Range("E48").Value = "'" + '  -2.55%  '
$ws.Range("D49").Value = "'" + '0.0986'
$ws.Range("E49").Value = "'" + '  -1.28%  '
$ws.Range("B50").Value = 'WOONetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D50").Value = "'" + '0.468'
$ws.Range("E50").Value = "'" + '  +7.23%  '
$ws.Range("D51").Value = "'" + '100.01'
$ws.Range("E51").Value = "'" + '  -4.24%  '
